# Fruta / hortaliza, semanal
#
# Inserts one new week's block of 4 rows (Especial / Primera / Segunda /
# Tercera) for date 44491 right before the existing row 646, pushing the
# remaining data down by 4 rows (old row 646 -> new row 650, ...,
# old row 750 -> new row 754). The new rows reuse the same static
# metadata (Mercado/Region/Tipo/Producto/Categoria/Variedad/Origen) as the
# surrounding "Piña" / "Caramelo" / "Ecuador" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 646:750 down by inserting 4 blank rows at row 646.
$ws.Rows("646:649").Insert()

# New data block (date 44491).
$newRows = @(
    @{ Row = 646; L = "Especial"; M = 25; N = 19000; O = 20000; P = 19600; Q = "$/caja 10 unidades"; S = 1960; T = 10 },
    @{ Row = 647; L = "Primera";  M = 30; N = 19000; O = 20000; P = 19500; Q = "$/caja 12 unidades"; S = 1625; T = 12 },
    @{ Row = 648; L = "Segunda";  M = 20; N = 19000; O = 20000; P = 19500; Q = "$/caja 14 unidades"; S = 1393; T = 14 },
    @{ Row = 649; L = "Tercera";  M = 30; N = 19000; O = 20000; P = 19500; Q = "$/caja 16 unidades"; S = 1219; T = 16 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = 44491
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100108
    $ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($row, 9).Value = 100108005
    $ws.Cells.Item($row, 10).Value = "Piña"
    $ws.Cells.Item($row, 11).Value = "Caramelo"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Ecuador"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
